# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.157.77"
Set-TextValue $ws.Range("E2") "  -2.88%  "
Set-TextValue $ws.Range("D3") "2.480.23"
Set-TextValue $ws.Range("E3") "  -3.17%  "
Set-TextValue $ws.Range("E4") "  -0.05%  "
Set-TextValue $ws.Range("D5") "564.28"
Set-TextValue $ws.Range("E5") "  -3.15%  "
Set-TextValue $ws.Range("D6") "163.20"
Set-TextValue $ws.Range("E6") "  -4.82%  "
Set-TextValue $ws.Range("E7") "  -0.06%  "
Set-TextValue $ws.Range("E8") "  -1.89%  "
Set-TextValue $ws.Range("D9") "2.478.93"
Set-TextValue $ws.Range("E9") "  -3.14%  "
Set-TextValue $ws.Range("E10") "  -5.59%  "
Set-TextValue $ws.Range("E11") "  -0.70%  "
Set-TextValue $ws.Range("E12") "  -2.57%  "
Set-TextValue $ws.Range("E13") "  -0.94%  "
Set-TextValue $ws.Range("D14") "2.937.32"
Set-TextValue $ws.Range("E14") "  -3.23%  "
Set-TextValue $ws.Range("D15") "69.038.98"
Set-TextValue $ws.Range("E15") "  -3.01%  "
Set-TextValue $ws.Range("E16") "  -3.20%  "
Set-TextValue $ws.Range("D18") "2.476.29"
Set-TextValue $ws.Range("E18") "  -4.23%  "
Set-TextValue $ws.Range("D19") "11.12"
Set-TextValue $ws.Range("E19") "  -4.47%  "
Set-TextValue $ws.Range("D20") "7.36"
Set-TextValue $ws.Range("E20") "  -7.70%  "
Set-TextValue $ws.Range("D21") "344.48"
Set-TextValue $ws.Range("E21") "  -3.65%  "
Set-TextValue $ws.Range("E22") "  -3.43%  "
Set-TextValue $ws.Range("E23") "  -7.52%  "
Set-TextValue $ws.Range("E24") "  -0.22%  "
Set-TextValue $ws.Range("D25") "69.29"
Set-TextValue $ws.Range("E25") "  -1.87%  "
Set-TextValue $ws.Range("D26") "3.86"
Set-TextValue $ws.Range("E26") "  -5.98%  "
Set-TextValue $ws.Range("D27") "2.607.80"
Set-TextValue $ws.Range("E27") "  -3.66%  "
Set-TextValue $ws.Range("D28") "8.62"
Set-TextValue $ws.Range("E28") "  -5.73%  "
Set-TextValue $ws.Range("D29") "0.999"
Set-TextValue $ws.Range("E29") "  -0.09%  "
Set-TextValue $ws.Range("E30") "  -6.35%  "
Set-TextValue $ws.Range("E31") "  -4.35%  "
Set-TextValue $ws.Range("D32") "439.55"
Set-TextValue $ws.Range("E32") "  -7.48%  "
Set-TextValue $ws.Range("B33") "FirstDigitalUSD"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D33") "0.999"
Set-TextValue $ws.Range("E33") "  -0.07%  "
Set-TextValue $ws.Range("B34") "Fetch.AI"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D34") "1.17"
Set-TextValue $ws.Range("E34") "  -8.72%  "
Set-TextValue $ws.Range("E35") "  -4.40%  "
Set-TextValue $ws.Range("D36") "155.21"
Set-TextValue $ws.Range("E36") "  -1.32%  "
Set-TextValue $ws.Range("E37") "  -5.28%  "
Set-TextValue $ws.Range("E38") "  -0.59%  "
Set-TextValue $ws.Range("D39") "18.06"
Set-TextValue $ws.Range("E39") "  -4.25%  "
Set-TextValue $ws.Range("E40") "  +0.01%  "
Set-TextValue $ws.Range("E41") "  -3.17%  "
Set-TextValue $ws.Range("D42") "4.55"
Set-TextValue $ws.Range("E42") "  -7.08%  "
Set-TextValue $ws.Range("E43") "  -3.88%  "
Set-TextValue $ws.Range("D44") "37.84"
Set-TextValue $ws.Range("E44") "  -2.36%  "
Set-TextValue $ws.Range("E45") "  -9.61%  "
Set-TextValue $ws.Range("E46") "  -9.37%  "
Set-TextValue $ws.Range("D47") "137.79"
Set-TextValue $ws.Range("E47") "  -5.67%  "
Set-TextValue $ws.Range("E48") "  -4.18%  "
Set-TextValue $ws.Range("E49") "  -5.88%  "
Set-TextValue $ws.Range("E50") "  -2.39%  "
Set-TextValue $ws.Range("D51") "0.569"
Set-TextValue $ws.Range("E51") "  -2.88%  "
